# Add an "intendedUse" field as a new row (row 7) on the active sheet,
# pushing the existing table (header + data, formerly rows 8-30) down by
# one row (now rows 9-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7 (shifts everything below down by one).
$ws.Rows.Item(7).Insert()

# Populate the new row with the intended-use label/value pair.
$ws.Cells.Item(7, 1).Value = "intendedUse"
$ws.Cells.Item(7, 2).Value = "Epi-validated outbreak"

# Match the row height used by the other label rows above it (15.75pt).
$ws.Rows.Item(7).RowHeight = 15.75
